# Updated symbol list on Tue Feb 14 22:12:47 UTC 2023 with GitHub Actions
#
# Refreshes the crypto-ranking snapshot on Sheet1 (coin order, prices,
# 1h volume %, and the "Hora" hour stamp) to match the latest scrape.
#
# All data cells in this sheet are stored as literal text (not numbers),
# including price/percentage/hour columns. A leading apostrophe is used
# when assigning values that look like numbers so Excel records them as
# text (quoted-number entry) instead of silently converting them to the
# Number type - this preserves the original text formatting (e.g. trailing
# zeros such as "296.66" or "0.00002100").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '296.66'
$ws.Range("E2").Value = "'" + '2.08%'
$ws.Range("G2").Value = "'" + '22'
$ws.Range("D3").Value = "'" + '41.61'
$ws.Range("E3").Value = "'" + '3.15%'
$ws.Range("G3").Value = "'" + '22'
$ws.Range("D4").Value = "'" + '5.039'
$ws.Range("E4").Value = "'" + '-0.11%'
$ws.Range("G4").Value = "'" + '22'
$ws.Range("D5").Value = "'" + '0.07556'
$ws.Range("E5").Value = "'" + '3.72%'
$ws.Range("G5").Value = "'" + '22'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").Value = "'" + '1.600'
$ws.Range("E6").Value = "'" + '2.36%'
$ws.Range("G6").Value = "'" + '22'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = "'" + '0.9286'
$ws.Range("E7").Value = "'" + '0.89%'
$ws.Range("G7").Value = "'" + '22'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = "'" + '2.410'
$ws.Range("E8").Value = "'" + '3.30%'
$ws.Range("G8").Value = "'" + '22'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'" + '0.1196'
$ws.Range("E9").Value = "'" + '3.42%'
$ws.Range("G9").Value = "'" + '22'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'" + '0.1842'
$ws.Range("E10").Value = "'" + '6.62%'
$ws.Range("G10").Value = "'" + '22'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'" + '0.08936'
$ws.Range("E11").Value = "'" + '3.73%'
$ws.Range("G11").Value = "'" + '22'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'" + '0.04036'
$ws.Range("E12").Value = "'" + '-3.57%'
$ws.Range("G12").Value = "'" + '22'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'" + '0.1050'
$ws.Range("E13").Value = "'" + '-0.25%'
$ws.Range("G13").Value = "'" + '22'
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").Value = "'" + '0.005981'
$ws.Range("E14").Value = "'" + '2.70%'
$ws.Range("G14").Value = "'" + '22'
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").Value = "'" + '3.349'
$ws.Range("E15").Value = "'" + '-1.50%'
$ws.Range("G15").Value = "'" + '22'
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").Value = "'" + '4.393'
$ws.Range("E16").Value = "'" + '2.54%'
$ws.Range("G16").Value = "'" + '22'
$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").Value = "'" + '0.3320'
$ws.Range("E17").Value = "'" + '1.30%'
$ws.Range("G17").Value = "'" + '22'
$ws.Range("B18").Value = 'MCDex'
$ws.Range("C18").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D18").Value = "'" + '7.919'
$ws.Range("E18").Value = "'" + '0.88%'
$ws.Range("G18").Value = "'" + '22'
$ws.Range("B19").Value = 'ProBitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D19").Value = "'" + '0.1419'
$ws.Range("E19").Value = "'" + '2.73%'
$ws.Range("G19").Value = "'" + '22'
$ws.Range("B20").Value = 'ZBToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D20").Value = "'" + '0.2997'
$ws.Range("E20").Value = "'" + '3.89%'
$ws.Range("G20").Value = "'" + '22'
$ws.Range("B21").Value = 'BitForexToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D21").Value = "'" + '0.001283'
$ws.Range("E21").Value = "'" + '1.40%'
$ws.Range("G21").Value = "'" + '22'
$ws.Range("D22").Value = "'" + '0.04058'
$ws.Range("E22").Value = "'" + '5.12%'
$ws.Range("G22").Value = "'" + '22'
$ws.Range("D23").Value = "'" + '0.001266'
$ws.Range("E23").Value = "'" + '-0.27%'
$ws.Range("G23").Value = "'" + '22'
$ws.Range("D24").Value = "'" + '0.003977'
$ws.Range("E24").Value = "'" + '4.38%'
$ws.Range("G24").Value = "'" + '22'
$ws.Range("D25").Value = "'" + '0.0001231'
$ws.Range("E25").Value = "'" + '-3.98%'
$ws.Range("G25").Value = "'" + '22'
$ws.Range("E26").Value = "'" + '-0.06%'
$ws.Range("G26").Value = "'" + '22'
$ws.Range("G27").Value = "'" + '22'
$ws.Range("G28").Value = "'" + '22'
$ws.Range("G29").Value = "'" + '22'
$ws.Range("G30").Value = "'" + '22'
$ws.Range("G31").Value = "'" + '22'
$ws.Range("G32").Value = "'" + '22'
$ws.Range("G33").Value = "'" + '22'
$ws.Range("G34").Value = "'" + '22'
$ws.Range("G35").Value = "'" + '22'
$ws.Range("G36").Value = "'" + '22'
$ws.Range("G37").Value = "'" + '22'
$ws.Range("D38").Value = "'" + '0.02410'
$ws.Range("E38").Value = "'" + '4.15%'
$ws.Range("G38").Value = "'" + '22'
$ws.Range("D39").Value = "'" + '0.05214'
$ws.Range("E39").Value = "'" + '5.34%'
$ws.Range("G39").Value = "'" + '22'
$ws.Range("D40").Value = "'" + '0.006393'
$ws.Range("E40").Value = "'" + '-3.69%'
$ws.Range("G40").Value = "'" + '22'
$ws.Range("D41").Value = "'" + '0.007782'
$ws.Range("E41").Value = "'" + '1.36%'
$ws.Range("G41").Value = "'" + '22'
$ws.Range("D42").Value = "'" + '0.1330'
$ws.Range("E42").Value = "'" + '4.58%'
$ws.Range("G42").Value = "'" + '22'
$ws.Range("D43").Value = "'" + '0.007545'
$ws.Range("E43").Value = "'" + '2.69%'
$ws.Range("G43").Value = "'" + '22'
$ws.Range("D44").Value = "'" + '0.007838'
$ws.Range("E44").Value = "'" + '10.92%'
$ws.Range("G44").Value = "'" + '22'
$ws.Range("D45").Value = "'" + '0.3210'
$ws.Range("E45").Value = "'" + '10.50%'
$ws.Range("G45").Value = "'" + '22'
$ws.Range("D46").Value = "'" + '0.00006787'
$ws.Range("E46").Value = "'" + '5.88%'
$ws.Range("G46").Value = "'" + '22'
$ws.Range("E47").Value = "'" + '-0.12%'
$ws.Range("G47").Value = "'" + '22'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = "'" + '0.004200'
$ws.Range("E48").Value = "'" + '-0.05%'
$ws.Range("G48").Value = "'" + '22'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = "'" + '0.04628'
$ws.Range("E49").Value = "'" + '166.44%'
$ws.Range("G49").Value = "'" + '22'
$ws.Range("D50").Value = "'" + '0.00002100'
$ws.Range("E50").Value = "'" + '-0.12%'
$ws.Range("G50").Value = "'" + '22'
$ws.Range("D51").Value = "'" + '0.0002000'
$ws.Range("E51").Value = "'" + '-0.12%'
$ws.Range("G51").Value = "'" + '22'
